$wb = $excel.ActiveWorkbook

# Fix the revenue-estimate formulas on the DCF sheet so they pull analyst
# estimates from RESEARCH instead of referencing stale/broken cells.
$dcf = $wb.Worksheets.Item("DCF")

$dcf.Range("Q22").Formula = "=IF(ISBLANK(U22),IF(OR(ISBLANK(RESEARCH!G15),RESEARCH!G13<>Q21),IF(OR(ISBLANK(RESEARCH!H15),RESEARCH!H13<>Q21),O22*(1+AVERAGE(`$G`$23:O23)),RESEARCH!H15/1000000),RESEARCH!G15/1000000),U22)"
$dcf.Range("R22").Formula = "=IF(ISBLANK(V22), IF(OR(ISBLANK(RESEARCH!H15),RESEARCH!H13<>R21), Q22*(1+AVERAGE(`$G`$23:O23,Q23)), RESEARCH!H15/1000000), V22)"
$dcf.Range("S22").Formula = "=IF(ISBLANK(W22), R22*(1+AVERAGE(`$G`$23:O23,Q23:R23)), W22)"

# Switch the active sheet/selection back to DCF (from RESEARCH).
$dcf.Activate()
$dcf.Range("AC13").Select()
